$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the column width for the new column BE (57) to match the existing 12-wide columns
$ws.Columns("BE").ColumnWidth = 11.16

# Header cell BE1: "2024/11/04" (kept as text, matching the style of BD1)
$ws.Range("BE1").NumberFormat = "@"
$ws.Range("BE1").Value = "2024/11/04"
$ws.Range("BD1").Copy()
$ws.Range("BE1").PasteSpecial(-4122)

$ws.Range("N2").Copy()
$ws.Range("BE2").PasteSpecial(-4122)
$ws.Range("BE2").Value = 126
$ws.Range("D2").Copy()
$ws.Range("BE3").PasteSpecial(-4122)
$ws.Range("BE3").Value = 115.7
$ws.Range("N2").Copy()
$ws.Range("BE4").PasteSpecial(-4122)
$ws.Range("BE4").Value = 127.3
$ws.Range("B2").Copy()
$ws.Range("BE5").PasteSpecial(-4122)
$ws.Range("BE5").Value = 143.8
$ws.Range("B2").Copy()
$ws.Range("BE6").PasteSpecial(-4122)
$ws.Range("BE6").Value = 156.5
$ws.Range("B2").Copy()
$ws.Range("BE7").PasteSpecial(-4122)
$ws.Range("BE7").Value = 167.4
$ws.Range("B2").Copy()
$ws.Range("BE8").PasteSpecial(-4122)
$ws.Range("BE8").Value = 171.9
$ws.Range("B2").Copy()
$ws.Range("BE9").PasteSpecial(-4122)
$ws.Range("BE9").Value = 149.7
$ws.Range("D2").Copy()
$ws.Range("BE10").PasteSpecial(-4122)
$ws.Range("BE10").Value = 113.1
$ws.Range("B2").Copy()
$ws.Range("BE11").PasteSpecial(-4122)
$ws.Range("BE11").Value = 161.3
$ws.Range("B2").Copy()
$ws.Range("BE12").PasteSpecial(-4122)
$ws.Range("BE12").Value = 191.7
$ws.Range("B2").Copy()
$ws.Range("BE13").PasteSpecial(-4122)
$ws.Range("BE13").Value = 212.4
$ws.Range("B2").Copy()
$ws.Range("BE14").PasteSpecial(-4122)
$ws.Range("BE14").Value = 166.1
$ws.Range("B2").Copy()
$ws.Range("BE15").PasteSpecial(-4122)
$ws.Range("BE15").Value = 150.9
$ws.Range("B2").Copy()
$ws.Range("BE16").PasteSpecial(-4122)
$ws.Range("BE16").Value = 156.1
$ws.Range("B2").Copy()
$ws.Range("BE17").PasteSpecial(-4122)
$ws.Range("BE17").Value = 152.5
$ws.Range("B2").Copy()
$ws.Range("BE18").PasteSpecial(-4122)
$ws.Range("BE18").Value = 184.1
$ws.Range("B2").Copy()
$ws.Range("BE19").PasteSpecial(-4122)
$ws.Range("BE19").Value = 212
$ws.Range("N2").Copy()
$ws.Range("BE20").PasteSpecial(-4122)
$ws.Range("BE20").Value = 130
$ws.Range("B2").Copy()
$ws.Range("BE21").PasteSpecial(-4122)
$ws.Range("BE21").Value = 168.7
$ws.Range("B2").Copy()
$ws.Range("BE22").PasteSpecial(-4122)
$ws.Range("BE22").Value = 154.3
$ws.Range("N2").Copy()
$ws.Range("BE23").PasteSpecial(-4122)
$ws.Range("BE23").Value = 133.1
$ws.Range("B2").Copy()
$ws.Range("BE24").PasteSpecial(-4122)
$ws.Range("BE24").Value = 163.9
$ws.Range("B2").Copy()
$ws.Range("BE25").PasteSpecial(-4122)
$ws.Range("BE25").Value = 187.7
$ws.Range("B2").Copy()
$ws.Range("BE26").PasteSpecial(-4122)
$ws.Range("BE26").Value = 147.3
$ws.Range("B2").Copy()
$ws.Range("BE27").PasteSpecial(-4122)
$ws.Range("BE27").Value = 187.5
$ws.Range("N2").Copy()
$ws.Range("BE28").PasteSpecial(-4122)
$ws.Range("BE28").Value = 132.2
$ws.Range("B2").Copy()
$ws.Range("BE29").PasteSpecial(-4122)
$ws.Range("BE29").Value = 154.3
$ws.Range("B2").Copy()
$ws.Range("BE30").PasteSpecial(-4122)
$ws.Range("BE30").Value = 150.3
$ws.Range("B2").Copy()
$ws.Range("BE31").PasteSpecial(-4122)
$ws.Range("BE31").Value = 190.3
$ws.Range("B2").Copy()
$ws.Range("BE32").PasteSpecial(-4122)
$ws.Range("BE32").Value = 176.7
$ws.Range("B2").Copy()
$ws.Range("BE33").PasteSpecial(-4122)
$ws.Range("BE33").Value = 195.9
$ws.Range("B2").Copy()
$ws.Range("BE34").PasteSpecial(-4122)
$ws.Range("BE34").Value = 142.3
$ws.Range("B2").Copy()
$ws.Range("BE35").PasteSpecial(-4122)
$ws.Range("BE35").Value = 156.5
$ws.Range("B2").Copy()
$ws.Range("BE36").PasteSpecial(-4122)
$ws.Range("BE36").Value = 141.8
$ws.Range("B2").Copy()
$ws.Range("BE37").PasteSpecial(-4122)
$ws.Range("BE37").Value = 149.7
$ws.Range("B2").Copy()
$ws.Range("BE38").PasteSpecial(-4122)
$ws.Range("BE38").Value = 154.5
$ws.Range("B2").Copy()
$ws.Range("BE39").PasteSpecial(-4122)
$ws.Range("BE39").Value = 155.8
$ws.Range("N2").Copy()
$ws.Range("BE40").PasteSpecial(-4122)
$ws.Range("BE40").Value = 136.4
$ws.Range("D2").Copy()
$ws.Range("BE41").PasteSpecial(-4122)
$ws.Range("BE41").Value = 119.4
$ws.Range("B2").Copy()
$ws.Range("BE42").PasteSpecial(-4122)
$ws.Range("BE42").Value = 142.7
$ws.Range("B2").Copy()
$ws.Range("BE43").PasteSpecial(-4122)
$ws.Range("BE43").Value = 157
$ws.Range("B2").Copy()
$ws.Range("BE44").PasteSpecial(-4122)
$ws.Range("BE44").Value = 143.4
$ws.Range("B2").Copy()
$ws.Range("BE45").PasteSpecial(-4122)
$ws.Range("BE45").Value = 142.1
$ws.Range("B2").Copy()
$ws.Range("BE46").PasteSpecial(-4122)
$ws.Range("BE46").Value = 395.4
$ws.Range("B2").Copy()
$ws.Range("BE47").PasteSpecial(-4122)
$ws.Range("BE47").Value = 151.1
$ws.Range("B2").Copy()
$ws.Range("BE48").PasteSpecial(-4122)
$ws.Range("BE48").Value = 170.3
$ws.Range("N2").Copy()
$ws.Range("BE49").PasteSpecial(-4122)
$ws.Range("BE49").Value = 136
$ws.Range("B2").Copy()
$ws.Range("BE50").PasteSpecial(-4122)
$ws.Range("BE50").Value = 146.9
$ws.Range("B2").Copy()
$ws.Range("BE51").PasteSpecial(-4122)
$ws.Range("BE51").Value = 209.4
$ws.Range("B2").Copy()
$ws.Range("BE52").PasteSpecial(-4122)
$ws.Range("BE52").Value = 182.9
$ws.Range("B2").Copy()
$ws.Range("BE53").PasteSpecial(-4122)
$ws.Range("BE53").Value = 175.1
